$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM-derived values per row/column, as described by the diff.
$updates = @{
    2 = @{ I=0.7104886604324562; J=0.7104886604324562; M=5.273684; N=15.821052; O=0.0510821201937383; P=0.0510821201937383; Q=14.100892300248; R=126.908030702232; S=0.03629326714849884; T=0.03629326714849884 }
    3 = @{ I=0.7104886604324562; J=0.7104886604324562; O=0.5598845502029881; P=0.5598845502029881; S=0.3977916240705493; T=0.3977916240705493 }
    4 = @{ I=0.7104886604324562; J=0.7104886604324562; M=32.95839133333334; N=98.87517400000002; O=0.3192425840231603; P=0.3192425840231604; Q=88.12487183167602; R=793.1238464850842; S=0.226818235875611; T=0.226818235875611 }
    5 = @{ I=0.7104886604324562; J=0.7104886604324562; M=7.205150000000001; N=21.61545; O=0.06979074558011317; P=0.06979074558011318; Q=19.2652885833; R=173.3875972497; S=0.04958553333779697; T=0.04958553333779698 }
    6 = @{ G=1.089534333333333; H=3.268603; I=0.2895113395675438; J=0.2895113395675438; M=5.273684; N=15.821052; O=0.0510821201937383; P=0.0510821201937383; Q=5.745859781150667; R=51.712738030356; S=0.01478885304523946; T=0.01478885304523946 }
    7 = @{ G=1.089534333333333; H=3.268603; I=0.2895113395675438; J=0.2895113395675438; O=0.5598845502029881; P=0.5598845502029881; Q=62.97738047868511; R=566.796424308166; S=0.1620929261324388; T=0.1620929261324388 }
    8 = @{ G=1.089534333333333; H=3.268603; I=0.2895113395675438; J=0.2895113395675438; M=32.95839133333334; N=98.87517400000002; O=0.3192425840231603; P=0.3192425840231604; Q=35.90929892910246; R=323.1836903619221; S=0.09242434814754931; T=0.09242434814754932 }
    9 = @{ G=1.089534333333333; H=3.268603; I=0.2895113395675438; J=0.2895113395675438; M=7.205150000000001; N=21.61545; O=0.06979074558011317; P=0.06979074558011318; Q=7.850258301816668; R=70.65232471635001; S=0.0202052122423162; T=0.02020521224231621 }
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $ws.Range("$col$row").Value = $updates[$row][$col]
    }
}
